$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, copy the "In Progress" look (currently on F8:F10) onto F12, which is
# moving from "Not Started" to "In Progress".
$ws.Range("F10").Copy()
$ws.Range("F12").PasteSpecial(-4122)
$ws.Range("F12").Value = "In Progress"

# Now flip the previously "In Progress" rows (Code View/Model/Controller) and
# the "Demo" row over to the "Completed" look (currently used by F3:F7).
$ws.Range("F3").Copy()
$ws.Range("F8:F11").PasteSpecial(-4122)
foreach ($r in 8..11) {
    $ws.Cells.Item($r, 6).Value = "Completed"
}

$excel.CutCopyMode = 0

# Move the active selection to reflect where the user ended up working.
$ws.Range("F9:F11").Select()
